# member_application.docx: bump the big centered year "2016" -> "2017"
# and move the "_GoBack" last-edit-position bookmark so it sits right
# after the new year text (this is what Word does automatically when a
# user's last edit was typing that "7"). All the form-field bookmarks
# (Check1, Check2, Text2.., Dropdown1, ...) keep their names but Word
# renumbers their w:id by +1 to make room for the new "_GoBack" bookmark,
# which always becomes id 0.

$d = $word.ActiveDocument

# --- 1. Change the digit run "6" -> "7" (turns "2016" into "2017") ---
# The big blue "20" / "1" / "6" runs live in paragraph 2; the lone "6"
# is the very last character of that paragraph, right before the
# paragraph mark.
$p = $d.Paragraphs(2)
$digitRange = $d.Range($p.Range.End - 2, $p.Range.End - 1)
$digitRange.Text = "7"

# --- 2. Re-seat the "_GoBack" bookmark right after the new "7" ---
# Collapsed (zero-length) ranges that sit exactly on a paragraph mark
# are mishandled when used directly with Bookmarks.Add, so temporarily
# insert a placeholder character after the "7", anchor the bookmark
# against that safe (non-edge) offset, then remove the placeholder.
# Re-adding a bookmark named "_GoBack" automatically replaces/removes
# whatever bookmark previously held that name (its old location, down
# near the "Text17" field), exactly matching Word's own behavior.
$p = $d.Paragraphs(2)
$placeholderPoint = $d.Range($p.Range.End - 1, $p.Range.End - 1)
$placeholderPoint.InsertAfter("Z")

$p = $d.Paragraphs(2)
$bmPos = $p.Range.End - 2
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

$p = $d.Paragraphs(2)
$placeholderRange = $d.Range($p.Range.End - 2, $p.Range.End - 1)
$placeholderRange.Delete()
